$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix number format of existing rows 60:61 (C:Z) so style 6 (numFmt 179) is no longer used ---
$ws.Range("C60:Z61").NumberFormat = "0.00_);[Red]\(0.00\)"

# --- Step 2: set default column style for C:Z to match (numFmt 177) ---
$ws.Range("C62:Z63").NumberFormat = "0.00_);[Red]\(0.00\)"

# --- Step 3: add new data rows 62 and 63 ---
# Row 62
$ws.Cells.Item(62, 1).Value = 45931
$ws.Cells.Item(62, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(62, 3).Value = 1013.4469999999997
$ws.Cells.Item(62, 4).Value = 1001.6139999999999
$ws.Cells.Item(62, 5).Value = 553.38499999999999
$ws.Cells.Item(62, 6).Value = 632.12200000000007
$ws.Cells.Item(62, 7).Value = 388.05899999999997
$ws.Cells.Item(62, 8).Value = 680.71
$ws.Cells.Item(62, 9).Value = 516.81200000000001
$ws.Cells.Item(62, 10).Value = 313.529
$ws.Cells.Item(62, 11).Value = 206.36699999999999
$ws.Cells.Item(62, 12).Value = 121.34700000000001
$ws.Cells.Item(62, 13).Value = 120.15000000000002
$ws.Cells.Item(62, 14).Value = 114.52999999999999
$ws.Cells.Item(62, 15).Value = 910.47899999999993
$ws.Cells.Item(62, 16).Value = 1230.4620000000002
$ws.Cells.Item(62, 17).Value = 656.68100000000004
$ws.Cells.Item(62, 18).Value = 483.54799999999994
$ws.Cells.Item(62, 19).Value = 312.60599999999999
$ws.Cells.Item(62, 20).Value = 397.10299999999995
$ws.Cells.Item(62, 21).Value = 166.01
$ws.Cells.Item(62, 22).Value = 175.94699999999997
$ws.Cells.Item(62, 23).Value = 195.047
$ws.Cells.Item(62, 24).Value = 100.42
$ws.Cells.Item(62, 25).Value = 78.039999999999992
$ws.Cells.Item(62, 26).Value = 0

# Row 63
$ws.Cells.Item(63, 1).Value = 45931
$ws.Cells.Item(63, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(63, 3).Value = 484.98399999999998
$ws.Cells.Item(63, 4).Value = 592.49200000000008
$ws.Cells.Item(63, 5).Value = 375.18700000000001
$ws.Cells.Item(63, 6).Value = 121.37
$ws.Cells.Item(63, 7).Value = 91.953000000000003
$ws.Cells.Item(63, 8).Value = 231.90299999999999
$ws.Cells.Item(63, 9).Value = 206.29
$ws.Cells.Item(63, 10).Value = 260.23100000000005
$ws.Cells.Item(63, 11).Value = 276.42499999999995
$ws.Cells.Item(63, 12).Value = 148.411
$ws.Cells.Item(63, 13).Value = 203.16099999999997
$ws.Cells.Item(63, 14).Value = 167.84099999999998
$ws.Cells.Item(63, 15).Value = 269.83699999999999
$ws.Cells.Item(63, 16).Value = 320.24700000000001
$ws.Cells.Item(63, 17).Value = 337.99700000000001
$ws.Cells.Item(63, 18).Value = 99.02600000000001
$ws.Cells.Item(63, 19).Value = 221.66399999999999
$ws.Cells.Item(63, 20).Value = 128.18799999999999
$ws.Cells.Item(63, 21).Value = 67.52
$ws.Cells.Item(63, 22).Value = 40.704999999999998
$ws.Cells.Item(63, 23).Value = 61.893000000000001
$ws.Cells.Item(63, 24).Value = 54.847000000000001
$ws.Cells.Item(63, 25).Value = 0
$ws.Cells.Item(63, 26).Value = 0

# --- Step 4: update selection ---
$ws.Range("K74").Select()